# "change speed to reach valid solution" — process time / 2, every agent go to TOP
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Process time (column D) is halved (rounded to nearest integer) for every
# data row (2-52): this is the "speed x2" / "process time / 2" change.
$newD = 11,22,34,45,56,67,78,89,100,111,122,134,145,156,167,`
        11,22,34,45,56,67,78,89,100,111,122,134,145,156,167,`
        11,22,34,45,56,67,78,89,100,111,122,134,145,156,167,`
        11,11,11,11,11,0

for ($i = 0; $i -lt $newD.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $newD[$i]
}

# Column E picked up an (empty) best-fit width entry from the user's session.
$ws.Columns("E:E").ColumnWidth = 8.786

# "every agent go to TOP": scroll the view back up and move the selection.
$ws.Range("G35").Select() | Out-Null
